$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$sm = $p.SlideMaster
$oldNmTheme = $nm.Theme
$oldSmTheme = $sm.Theme
$nm.Theme = $oldSmTheme
$sm.Theme = $oldNmTheme
Write-Output "done"
